$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update input values on the "Inputs" side of the sheet ---
$ws.Range("B9").Value = 29.6131    # PV cost, installed ($/kW)
$ws.Range("B10").Value = 1.69242   # Battery cost ($/kWh)
$ws.Range("B11").Value = 1.35394   # Battery inverter cost ($/kW)
$ws.Range("B13").Value = 2050      # PV O&M ($/kW/yr)
$ws.Range("B16").Value = 12.6      # Battery replacement cost
$ws.Range("B19").Value = 4042      # Year 1 energy savings
$ws.Range("B20").Value = 2818      # Year 1 demand savings
$ws.Range("B21").Value = 20        # analysis_period (years) - was 25

# --- Shrink the cash-flow table from a 25 year horizon (cols B:AA) to the
#     new 20 year horizon (cols B:V) by deleting the now-unused columns
#     W:AA. Those columns hold no data outside the cash-flow block, so
#     this only removes the trailing year-21..25 cells and naturally
#     recomputes the sheet's used range / row spans. ---
$ws.Range("W1:AA39").Delete()

# --- Column A narrows (its label column got tighter); column B (and the
#     rest of the yearly cash-flow block) keeps its existing width. ---
$ws.Columns.Item(1).ColumnWidth = 20.8333333333333
